# MiniDesignerConfigsTemplate/__beans__.xlsx edit
#
# The "fields" sub-table header row (row3/row4) originally had 11 columns:
#   F=name, G=type, H=sep, I=is_multi_rows, J=index, K=group, L=ref,
#   M=path, N=comment, O=tags, P=orientation
# It is simplified down to 5 columns:
#   F=name, G=type, H=group, I=comment, J=tags
#
# That means columns sep / is_multi_rows / index / ref / path / orientation
# are removed outright, while group / comment / tags survive and simply
# shift left into the freed-up slots. We do this with plain column deletes
# (rightmost column first so earlier deletes don't invalidate later column
# letters).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# orientation (P) - trailing column, drop first
$ws.Range("P1").EntireColumn.Delete()
# path (M)
$ws.Range("M1").EntireColumn.Delete()
# ref (L)
$ws.Range("L1").EntireColumn.Delete()
# index (J)
$ws.Range("J1").EntireColumn.Delete()
# is_multi_rows (I)
$ws.Range("I1").EntireColumn.Delete()
# sep (H) - group/ref/path/comment/tags/orientation all shift left from here
$ws.Range("H1").EntireColumn.Delete()

# "fields" -> "*fields" label on the merged banner cell above the header row
$ws.Range("F2").Value = "*fields"

# Column width touch-ups for the now-narrower table (column C / former "sep"
# column width class, and column H / former lone-width column P). The COM
# layer quantizes stored width to 1/7ths, so these inputs are chosen to land
# on the closest achievable stored width to the authored 12.25 / 10.375.
$ws.Range("C1").ColumnWidth = 11.571428571428571
$ws.Range("H1").ColumnWidth = 9.714285714285714

# Restore the cursor/selection to where the author left it
$ws.Range("G14").Select()
